# Updated cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be treated as text so numeric-looking strings like
    # "211.08" or "7.00" are not coerced into numbers (which would lose the
    # trailing zero / precision of the original formatted price string).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.609.31"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.597.48"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.17%  "

# Row 5 - BNB
Set-TextValue "D5" "211.08"
$ws.Range("E5").Value = "  -0.20%  "

# Row 6 - XRP
Set-TextValue "D6" "0.514"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.0618"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.246"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Solana
Set-TextValue "D10" "19.47"
$ws.Range("E10").Value = "  -1.21%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0841"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.821.89"

# Row 13 - WrappedEther
Set-TextValue "D13" "1.605.12"
$ws.Range("E13").Value = "  +0.60%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.02"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.39%  "

# Row 16 - Litecoin
Set-TextValue "D16" "64.74"
$ws.Range("E16").Value = "  -0.14%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "26.609.40"
$ws.Range("E17").Value = "  -0.17%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.89%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.12%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "208.10"
$ws.Range("E20").Value = "  -0.67%  "

# Row 21 - Chainlink
Set-TextValue "D21" "7.00"
$ws.Range("E21").Value = "  +4.29%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.34%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -0.52%  "

# Row 24 - Avalanche
Set-TextValue "D24" "8.89"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25 - Monero
Set-TextValue "D25" "145.38"
$ws.Range("E25").Value = "  -0.95%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.20%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.55%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.53%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.26"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +1.56%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.15"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.22"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "2.94"
$ws.Range("E33").Value = "  +0.64%  "

# Row 34 - Maker
Set-TextValue "D34" "1.280.50"
$ws.Range("E34").Value = "  -0.99%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "0.615"
$ws.Range("E35").Value = "  -7.54%  "

# Row 36 - HuobiToken
Set-TextValue "D36" "2.46"
$ws.Range("E36").Value = "  +0.82%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +0.35%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.77%  "

# Row 39 - ARBITRUM
Set-TextValue "D39" "0.837"
$ws.Range("E39").Value = "  +0.15%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +20.82%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.47"
$ws.Range("E41").Value = "  +2.14%  "

# Row 42 - MXToken
Set-TextValue "D42" "2.20"
$ws.Range("E42").Value = "  +0.21%  "

# Row 43 & 44 - Aave / TrustWalletToken swap places
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D43" "0.784"
$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "64.15"
$ws.Range("E44").Value = "  +0.85%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.734.06"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46 - Quant
Set-TextValue "D46" "90.17"
$ws.Range("E46").Value = "  +0.36%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -0.74%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +3.88%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.82%  "

# Row 50 - USDD
$ws.Range("E50").Value = "  +0.10%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "7.44"
$ws.Range("E51").Value = "  -0.76%  "
